$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 8540.925999999999
$ws.Range("J17").Value = 8540.925999999999
$ws.Range("L17").Value = 25622.778
$ws.Range("N17").Value = -25958.778

# Row 80
$ws.Range("H80").Value = 8050.0835
$ws.Range("I80").Value = 893.2
$ws.Range("K80").Value = 2679.6
$ws.Range("M80").Value = -1681.6

# Row 83
$ws.Range("H83").Value = 8050.0835
$ws.Range("I83").Value = 893.2
$ws.Range("K83").Value = 8038.8
$ws.Range("M83").Value = -3046.8

# Row 86
$ws.Range("H86").Value = 4589.48
$ws.Range("I86").Value = 6419.1665
$ws.Range("J86").Value = 2900.5386
$ws.Range("K86").Value = 6419.1665
$ws.Range("L86").Value = 2900.5386
$ws.Range("M86").Value = -5296.1665
$ws.Range("N86").Value = -5146.5386

# Row 89
$ws.Range("H89").Value = 4589.48
$ws.Range("I89").Value = 6419.1665
$ws.Range("J89").Value = 2900.5386
$ws.Range("K89").Value = 32095.8325
$ws.Range("L89").Value = 14502.693
$ws.Range("M89").Value = -26479.8325
$ws.Range("N89").Value = -25734.693

# Row 138
$ws.Range("H138").Value = 2930.3022
$ws.Range("I138").Value = 1875.4286
$ws.Range("J138").Value = 3937.2273
$ws.Range("K138").Value = 5626.2858
$ws.Range("L138").Value = 11811.6819
$ws.Range("M138").Value = -486.2857999999997
$ws.Range("N138").Value = -22091.6819

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1917462.4
$ws.Range("I2").Value = 2647022.8
$ws.Range("K2").Value = 2647022.8
$ws.Range("M2").Value = -2646909.8

# Row 31
$ws.Range("H31").Value = 1900
$ws.Range("I31").Value = 1900
$ws.Range("K31").Value = 1900
$ws.Range("M31").Value = -1606

# Row 32
$ws.Range("H32").Value = 11052.521
$ws.Range("I32").Value = 6636.923
$ws.Range("K32").Value = 6636.923
$ws.Range("M32").Value = -6349.923

# Row 74
$ws.Range("H74").Value = 20797.633
$ws.Range("I74").Value = 2039.0714
$ws.Range("J74").Value = 133349
$ws.Range("K74").Value = 2039.0714
$ws.Range("L74").Value = 133349
$ws.Range("M74").Value = -1165.0714
$ws.Range("N74").Value = -135097

# Row 77
$ws.Range("H77").Value = 20797.633
$ws.Range("I77").Value = 2039.0714
$ws.Range("J77").Value = 133349
$ws.Range("K77").Value = 10195.357
$ws.Range("L77").Value = 666745
$ws.Range("M77").Value = -5827.357
$ws.Range("N77").Value = -675481

# Row 97
$ws.Range("H97").Value = 1291639.8
$ws.Range("I97").Value = 1459414.1
$ws.Range("J97").Value = 5370
$ws.Range("K97").Value = 1459414.1
$ws.Range("L97").Value = 5370
$ws.Range("M97").Value = -1458918.1
$ws.Range("N97").Value = -6362

# Row 109
$ws.Range("H109").Value = 65125.668
$ws.Range("J109").Value = 65125.668
$ws.Range("L109").Value = 65125.668
$ws.Range("N109").Value = -67899.66800000001

# Row 110
$ws.Range("H110").Value = 1323513.2
$ws.Range("I110").Value = 1323513.2
$ws.Range("K110").Value = 1323513.2
$ws.Range("M110").Value = -1321468.2

# Row 116
$ws.Range("H116").Value = 1917462.4
$ws.Range("I116").Value = 2647022.8
$ws.Range("K116").Value = 2647022.8
$ws.Range("M116").Value = -2644728.8

# Row 132
$ws.Range("H132").Value = 26843.277
$ws.Range("I132").Value = 1343.72
$ws.Range("K132").Value = 4031.16
$ws.Range("M132").Value = -1501.16

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1917462.4
$ws.Range("I3").Value = 2647022.8
$ws.Range("K3").Value = 2647022.8
$ws.Range("M3").Value = -2646908.8

# Row 24
$ws.Range("H24").Value = 10016
$ws.Range("I24").Value = 10016
$ws.Range("K24").Value = 10016
$ws.Range("M24").Value = -9781

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 17296.18
$ws.Range("J31").Value = 25800.07
$ws.Range("L31").Value = 25800.07
$ws.Range("N31").Value = -26390.07

# Row 34
$ws.Range("H34").Value = 17296.18
$ws.Range("J34").Value = 25800.07
$ws.Range("L34").Value = 25800.07
$ws.Range("N34").Value = -26204.07

# Row 134
$ws.Range("H134").Value = 3934.2666
$ws.Range("I134").Value = 1499.4
$ws.Range("K134").Value = 4498.200000000001
$ws.Range("M134").Value = -1963.200000000001

# Row 141
$ws.Range("H141").Value = 71196.88
$ws.Range("J141").Value = 78933.266
$ws.Range("L141").Value = 78933.266
$ws.Range("N141").Value = -89293.266

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 14885012
$ws.Range("I131").Value = 6946561
$ws.Range("J131").Value = 20838850
$ws.Range("K131").Value = 20839683
$ws.Range("L131").Value = 62516550
$ws.Range("M131").Value = -20834643
$ws.Range("N131").Value = -62526630

# Row 134
$ws.Range("H134").Value = 1388.9231
$ws.Range("I134").Value = 1343.5834
$ws.Range("K134").Value = 4030.7502
$ws.Range("M134").Value = 1039.2498

# Row 139
$ws.Range("H139").Value = 1918.5454
$ws.Range("I139").Value = 1851.1904
$ws.Range("J139").Value = 3333
$ws.Range("K139").Value = 5553.5712
$ws.Range("L139").Value = 9999
$ws.Range("M139").Value = -413.5712000000003
$ws.Range("N139").Value = -20279

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 1419520
$ws.Range("I80").Value = 2750520.8
$ws.Range("K80").Value = 2750520.8
$ws.Range("M80").Value = -2749522.8

# Row 83
$ws.Range("H83").Value = 1419520
$ws.Range("I83").Value = 2750520.8
$ws.Range("K83").Value = 13752604
$ws.Range("M83").Value = -13747612

# Row 113
$ws.Range("H113").Value = 6546667.5
$ws.Range("I113").Value = 10394972
$ws.Range("J113").Value = 4550
$ws.Range("K113").Value = 10394972
$ws.Range("L113").Value = 4550
$ws.Range("M113").Value = -10392802
$ws.Range("N113").Value = -8890

# Row 132
$ws.Range("H132").Value = 2550.418
$ws.Range("I132").Value = 2263.6597
$ws.Range("J132").Value = 3224.3
$ws.Range("K132").Value = 6790.9791
$ws.Range("L132").Value = 9672.900000000001
$ws.Range("M132").Value = -4260.9791
$ws.Range("N132").Value = -14732.9

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 114661.375
$ws.Range("I22").Value = 445044
$ws.Range("K22").Value = 445044
$ws.Range("M22").Value = -444749

# Row 26
$ws.Range("H26").Value = 5633.3335
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1705

# Row 27
$ws.Range("H27").Value = 114661.375
$ws.Range("I27").Value = 445044
$ws.Range("K27").Value = 445044
$ws.Range("M27").Value = -444937

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 31
$ws.Range("H31").Value = 20999.666
$ws.Range("J31").Value = 20999.666
$ws.Range("L31").Value = 20999.666
$ws.Range("N31").Value = -21695.666

# Row 41
$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 3000
$ws.Range("L41").Value = 3000
$ws.Range("N41").Value = -3780
